# Update cryptocurrency price/volume data per Aug 26 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so values such as
# "26.069.47", "1.001" or "16.00" keep their exact literal formatting
# instead of being auto-coerced into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.069.47"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.651.86"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "217.35"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.5277"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.2596"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "0.06314"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "20.32"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").Value = "0.07788"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.652.28"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "1.878.78"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "0.5476"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "0.0₅8183"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "26.073.77"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "4.573"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "190.75"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "6.013"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "143.61"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "0.1231"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "16.00"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").Value = "1.433"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "0.05798"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "3.541"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "3.258"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "1.592"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Value = "2.795"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "2.411"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "0.9420"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").Value = "0.5749"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "104.89"
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("D41").Value = "0.8480"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "5.709"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("D44").Value = "1.026.88"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("D45").Value = "1.795.54"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "57.14"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "0.4328"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.845"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05140"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "1.447"
$ws.Range("E51").Value = "  -1.41%  "
